$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h change (E) columns for the cryptos list.
# A leading apostrophe forces the numeric-looking strings to be stored as text
# (matching the source data's inline-string cell type), then ClearFormats()
# removes the quote-prefix formatting Excel applies so the cell's style is
# left untouched.

$ws.Range('D2').Value = '''67.810.90'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '''  -0.98%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = '''3.737.18'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '''  -3.05%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = '''  -0.06%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = '''593.87'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '''  -1.30%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = '''166.14'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '''  -3.70%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = '''3.734.73'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '''  -3.08%  '
$ws.Range('E7').ClearFormats()
$ws.Range('E8').Value = '''  -0.07%  '
$ws.Range('E8').ClearFormats()
$ws.Range('E9').Value = '''  -2.13%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = '''0.160'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '''  -4.66%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = '''6.43'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '''  -1.82%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = '''0.446'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '''  -3.43%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = '''0.0000266'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '''  -6.17%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = '''35.98'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '''  -3.23%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = '''4.370.91'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '''  -2.91%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = '''3.736.10'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '''  -2.59%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = '''67.712.73'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '''  -1.17%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = '''18.33'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '''  -0.34%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = '''7.04'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '''  -5.67%  '
$ws.Range('E19').ClearFormats()
$ws.Range('E20').Value = '''  -0.45%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = '''10.45'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '''  -4.13%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = '''463.68'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '''  -1.57%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = '''0.702'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '''  -4.47%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = '''82.87'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '''  -1.09%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = '''0.0000137'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '''  -14.18%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = '''2.19'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '''  -3.95%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = '''11.99'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '''  -1.94%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = '''10.29'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '''  -2.88%  '
$ws.Range('E28').ClearFormats()
$ws.Range('E29').Value = '''  +0.02%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = '''2.88'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '''  -2.12%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = '''3.887.71'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '''  -2.98%  '
$ws.Range('E31').ClearFormats()
$ws.Range('D32').Value = '''7.39'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '''  -5.19%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = '''29.88'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '''  -4.18%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = '''2.18'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '''  -5.72%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').Value = '''9.04'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '''  -4.06%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').Value = '''3.689.14'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '''  -3.36%  '
$ws.Range('E36').ClearFormats()
$ws.Range('E37').Value = '''  -3.44%  '
$ws.Range('E37').ClearFormats()
$ws.Range('D38').Value = '''3.55'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '''  -9.54%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = '''0.137'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '''  -2.16%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = '''0.992'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '''  -2.93%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = '''5.73'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '''  -4.14%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').Value = '''1.00'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '''  +0.01%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D44').Value = '''0.306'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '''  -4.51%  '
$ws.Range('E44').ClearFormats()
$ws.Range('D45').Value = '''8.53'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '''  -2.69%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').Value = '''1.91'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '''  -4.30%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = '''396.83'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '''  -5.68%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = '''45.10'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '''  -3.61%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = '''143.91'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '''  +1.20%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = '''38.85'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '''  +1.20%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').Value = '''0.0347'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '''  -3.82%  '
$ws.Range('E51').ClearFormats()
